# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" quarterly sheet (right after "总计" and before the
# former first quarter sheet "2022-Q3"), populates it with the new quarter's
# fund-holding figures, and refreshes the "总计" (totals) summary sheet so its
# quarter labels / values shift down one row and a new trailing row for
# "2020-Q4" appears.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to be stored as TEXT
# (the fund sheets keep their numeric-looking figures as text, matching the
# source data), without leaving a lingering custom number-format style on
# the cell.
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet by duplicating the "2022-Q3" sheet
#    (same headers / fund universe / styles), then overwrite just the
#    figures that differ for the new quarter.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)            # new copy is placed immediately before 2022-Q3
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

Set-TextValue $q4.Range("D2") "10.01"
Set-TextValue $q4.Range("E2") "93.81"
Set-TextValue $q4.Range("F2") "1.07"
Set-TextValue $q4.Range("G2") "0.1071"
$q4.Range("H2").Value = 5

Set-TextValue $q4.Range("D3") "0.94"
Set-TextValue $q4.Range("E3") "97.66"
Set-TextValue $q4.Range("F3") "1.12"
Set-TextValue $q4.Range("G3") "0.0105"
$q4.Range("H3").Value = 5

# ---------------------------------------------------------------------------
# 2) Refresh the "总计" summary sheet: each existing quarter's figures move
#    down one row (column A's running index 0..5 is untouched - it is just a
#    positional counter), the new 2022-Q4 figures take row 2, and a new row 7
#    appears for what used to be the last row's quarter, 2020-Q4.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 7 is brand new - give it the same look (borders/bold) as the other
# index cells in column A before filling it in.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 2
$total.Range("D7").Value = 0.26

$total.Range("B6").Value = "2021-Q4"
$total.Range("D6").Value = 0.22

$total.Range("B5").Value = "2022-Q1"
$total.Range("D5").Value = 0.16

$total.Range("B4").Value = "2022-Q2"
$total.Range("D4").Value = 0.17

$total.Range("B3").Value = "2022-Q3"
$total.Range("D3").Value = 0.15

$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 0.12

$total.Select()
$total.Range("A1").Select()
